{"js": "async function replaceText(body, findText, newText) {\n  const results = body.search(findText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n\nconst body = context.document.body;\n\n// Title (appears twice: H1 heading + bold paragraph near the end) - replace both occurrences\nawait replaceText(\n  body,\n  \"Play Dr. Magoo's Adventure Slot for Free - Pros and Cons\",\n  \"Play Dr. Magoo's Adventure for Free - Exciting Slot Game\"\n);\n\n// \"What we like\" bullet list - text is shuffled/replaced in place; process in an order that\n// never lets an already-written replacement collide with a still-pending \"find\" target.\nawait replaceText(\n  body,\n  \"Excellent RTP value of 96.33%\",\n  \"High-resolution graphics and immersive tribal sound\"\n);\nawait replaceText(\n  body,\n  \"High-resolution graphics with attention to detail\",\n  \"Customizable gaming experience with various options in the internal menu\"\n);\nawait replaceText(\n  body,\n  \"Bonus game with higher value symbols and additional free spins\",\n  \"Excellent RTP value of 96.33%\"\n);\nawait replaceText(\n  body,\n  \"Wild symbol expands to fill entire reel creating numerous winning combos\",\n  \"Good balance between spins made and winning combinations\"\n);\n\n// \"What we don't like\" bullet list\nawait replaceText(\n  body,\n  \"Limited number of paylines at 20\",\n  \"Limited bonus game with restricted access to higher-value symbols\"\n);\nawait replaceText(\n  body,\n  \"No progressive jackpot feature\",\n  \"Limited availability of free spins in the base game\"\n);\n\n// Closing italic summary paragraph\nawait replaceText(\n  body,\n  \"Read our review on Dr. Magoo's Adventure, a free online slot game featuring high-resolution graphics, bonus games, free spins, and excellent RTP value.\",\n  \"Play Dr. Magoo's Adventure for free and enjoy high-resolution graphics, bonus features, and excellent RTP value.\"\n);\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# Locate the first occurrence of $findText at or after character position $searchStart\n# and return a fresh Range spanning exactly that match (Find itself, without a\n# Replacement string, leaves straight apostrophes alone - unlike Find.Execute with a\n# Replacement argument, which this runtime \"smart quotes\").\nfunction Find-Range($findText, $searchStart) {\n    $rng = $d.Range($searchStart, $d.Content.End)\n    $rng.Find.ClearFormatting()\n    $rng.Find.Forward = $true\n    $rng.Find.Wrap = 0\n    $rng.Find.Text = $findText\n    $found = $rng.Find.Execute()\n    if ($found) {\n        return $d.Range($rng.Start, $rng.End)\n    }\n    return $null\n}\n\n# Replace the first occurrence of $findText found at/after $searchStart with $newText.\n# Returns the character position right after the inserted text, so callers can keep\n# scanning forward (needed for \"Play Dr. Magoo's Adventure...\" which occurs twice).\nfunction Replace-FirstFrom($findText, $newText, $searchStart) {\n    $target = Find-Range $findText $searchStart\n    $target.Text = $newText\n    return $target.Start + $newText.Length\n}\n\n# Title (appears twice: the H1 heading and the bold paragraph near the end)\n$pos = Replace-FirstFrom \"Play Dr. Magoo's Adventure Slot for Free - Pros and Cons\" \"Play Dr. Magoo's Adventure for Free - Exciting Slot Game\" 0\n$pos = Replace-FirstFrom \"Play Dr. Magoo's Adventure Slot for Free - Pros and Cons\" \"Play Dr. Magoo's Adventure for Free - Exciting Slot Game\" $pos\n\n# \"What we like\" bullet list (4 items) - content reshuffled/replaced in place\n[void](Replace-FirstFrom \"Excellent RTP value of 96.33%\" \"High-resolution graphics and immersive tribal sound\" 0)\n[void](Replace-FirstFrom \"High-resolution graphics with attention to detail\" \"Customizable gaming experience with various options in the internal menu\" 0)\n[void](Replace-FirstFrom \"Bonus game with higher value symbols and additional free spins\" \"Excellent RTP value of 96.33%\" 0)\n[void](Replace-FirstFrom \"Wild symbol expands to fill entire reel creating numerous winning combos\" \"Good balance between spins made and winning combinations\" 0)\n\n# \"What we don't like\" bullet list (2 items)\n[void](Replace-FirstFrom \"Limited number of paylines at 20\" \"Limited bonus game with restricted access to higher-value symbols\" 0)\n[void](Replace-FirstFrom \"No progressive jackpot feature\" \"Limited availability of free spins in the base game\" 0)\n\n# Closing italic summary paragraph\n[void](Replace-FirstFrom \"Read our review on Dr. Magoo's Adventure, a free online slot game featuring high-resolution graphics, bonus games, free spins, and excellent RTP value.\" \"Play Dr. Magoo's Adventure for free and enjoy high-resolution graphics, bonus features, and excellent RTP value.\" 0)\n"}
